$d = $word.ActiveDocument

# 1. Author name typo fix: "Danie" -> "Daniel"
$d.Content.Find.Execute("Danie Easterling", $true, $false, $false, $false, $false, $true, 1, $false, "Daniel Easterling", 2) | Out-Null

# 2. "using the intranets of ... coffehouses," -> "using the Intranet of ... coffeehouses,"
$d.Content.Find.Execute("using the intranets of various computer labs, companies, coffehouses,", $true, $false, $false, $false, $false, $true, 1, $false, "using the Intranet of various computer labs, companies, coffeehouses,", 2) | Out-Null

# 3. "wherever an intranet may be." -> "wherever an Intranet may be."
$d.Content.Find.Execute("wherever an intranet may be.", $true, $false, $false, $false, $false, $true, 1, $false, "wherever an Intranet may be.", 2) | Out-Null

# 4. "confines of the intranet on which Bitter" -> "confines of the Intranet on which Bitter"
$d.Content.Find.Execute("confines of the intranet on which Bitter", $true, $false, $false, $false, $false, $true, 1, $false, "confines of the Intranet on which Bitter", 2) | Out-Null

# 5. "case-write up." -> "case write-up."
$d.Content.Find.Execute("functions using a case-write up.", $true, $false, $false, $false, $false, $true, 1, $false, "functions using a case write-up.", 2) | Out-Null

# 6. "sys-admin" -> "system administrator"
$d.Content.Find.Execute("The client of Bitter is the sys-admin on the up-and-up.", $true, $false, $false, $false, $false, $true, 1, $false, "The client of Bitter is the system administrator on the up-and-up.", 2) | Out-Null

# 7. "history which means" -> "history; this means"
$d.Content.Find.Execute("history which means", $true, $false, $false, $false, $false, $true, 1, $false, "history; this means", 2) | Out-Null

# 8. Split the paragraph: "... today." + " Bitter will enable ... writer's messages."
#    becomes two separate paragraphs, the first ending in " today. " (note trailing space) and a
#    second new paragraph (with a leading tab) holding the "Bitter will enable..." text.
$p = $d.Paragraphs.Item(9)
$full = $p.Range.Text
$marker = " Bitter will enable the user to create an account, and login if their credentials are correct. Given the user logged in correctly, create a profile, and the user can post messages to another user, and read messages in a conversation with another user. An external user to the conversation can opt to " + [char]8220 + "subscribe" + [char]8221 + " to the users conversation or message history; this means the external user receives a read-only log of the writer's messages."
$idx = $full.IndexOf($marker)
if ($idx -lt 0) {
    throw "marker text for paragraph split not found"
}
$rngStart = $p.Range.Start + $idx
$rngEnd = $rngStart + $marker.Length
$rng = $d.Range($rngStart, $rngEnd)
$rng.Text = " "
$rng.InsertParagraphAfter() | Out-Null

$p2 = $d.Paragraphs.Item(10)
$p2.Range.InsertAfter("`tBitter will enable the user to create an account, and login if their credentials are correct. Given the user logged in correctly, create a profile, and the user can post messages to another user, and read messages in a conversation with another user. An external user to the conversation can opt to " + [char]8220 + "subscribe" + [char]8221 + " to the users conversation or message history; this means the external user receives a read-only log of the writer's messages.")
